$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Remove the obsolete bottom rows (old row 5 "corev_rand_interrupt_nested"
# content plus the two blank rows 6 and 7). The two surviving data rows
# (3 and 4) get entirely new content below.
# ------------------------------------------------------------------
$ws.Rows("5:7").Delete()

# ------------------------------------------------------------------
# Row 3: "corev_directed_for_hwloop_covg_test"
# ------------------------------------------------------------------
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = ""
$ws.Range("D3").Value = "Instruction number in HWLOOP body exceeding maximum number of instruction when using immediate version of HWLOOP instructions (cv.counti, cv.setupi).`nThis issue seen rarely, maybe 1 out of 2~3 regressions. Too risky to fix it now as it might affect other tests. "
$ws.Range("E3").Value = "Test generation corner case"
$ws.Range("F3").Value = "Open (no fix)"
$ws.Range("G3").Value = "Random"
$ws.Range("H3").Value = "corev_directed_for_hwloop_covg_test"
$ws.Range("I3").Value = "test_program/corev_directed_for_hwloop_covg_test_550986902.S:826:(.text+0xbb4): relocation truncated to fit: R_RISCV_CVPCREL_UI12 against ``hwloop1_end_stream1_id0'"

# ------------------------------------------------------------------
# Row 4: "corev_rand_fp_instr_debug_test_with_int_debug_trigger_and_single_step"
# ------------------------------------------------------------------
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = ""
$ws.Range("D4").Value = "Load in Virtual Peripheral area`nx7 set to 0x0fffffce and going to 0x10000002 with some load address offset"
$ws.Range("E4").Value = "Test generation corner case"
$ws.Range("F4").Value = "Open (no fix)"
$ws.Range("G4").Value = "Random"
$ws.Range("H4").Value = 12
$ws.Range("I4").Value = "# (IDV) Instruction executed prior to mismatch '0x1a110d1e(debug_end+48): 02c3a587 flw     f11,44(x7)'"

# ------------------------------------------------------------------
# Re-normalise cell formatting.
#
# A handful of cells in the source workbook carried stray direct
# formatting (a date number format left over from a copy/paste, and a
# "treat as text" quote-prefix on the numeric S/N column) that Excel
# cleaned up when the sheet was refreshed. Re-apply plain formats by
# copying from a cell that already has the desired combination, which
# keeps the existing style table entries instead of inventing new ones.
# ------------------------------------------------------------------

# Plain bordered/aligned cells (S/N, ticket#, Category, Random/Custom,
# failing-test columns) -> copy the clean format already used by H4.
# (PasteSpecial only honours the first area of a multi-area Union, so
# each destination block is pasted individually.)
foreach ($addr in @("B3:C4", "E3:E4", "G3:G4", "H3")) {
    $ws.Range("H4").Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
}

# Wrapped-text bordered cells (Known issue / error message columns) ->
# copy the clean format already used by D4 (still has the original,
# untouched wrap style at this point).
foreach ($addr in @("D3", "I3:I4")) {
    $ws.Range("D4").Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
}

# Status column keeps the "Neutral" look (tan text on pale-yellow fill)
# but without the stray date format / word-wrap -> start from a plain
# bordered cell, then restore the Neutral colours directly.
$ws.Range("H4").Copy() | Out-Null
$ws.Range("F3:F4").PasteSpecial(-4122) | Out-Null
$ws.Range("F3:F4").Font.Color = 22428
$ws.Range("F3:F4").Interior.Color = 10284031

$excel.CutCopyMode = 0

$ws.Rows("3").RowHeight = 57.6
$ws.Rows("4").RowHeight = 28.8

# ------------------------------------------------------------------
# Sheet-level view tweaks.
# ------------------------------------------------------------------
$excel.ActiveWindow.Zoom = 90
$ws.Range("E20").Select() | Out-Null

# ------------------------------------------------------------------
# Column widths (content-driven change: column E now holds the longer
# "Test generation corner case" text instead of "Imperas"/"verif").
# ------------------------------------------------------------------
$ws.Columns("E").ColumnWidth = 25.77734375
